# CRM-1980 Add all upcountry, prepaid, postpaid, invoices and contract
# related fields in the Partner excel file which is downloaded from our Panel.
#
# Adds 14 new header/placeholder column pairs (Q..AD) to row 1 / row 2 of
# Sheet1, widens a few columns, sets the printed page to portrait, and
# restores the last-used selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New field definitions, in the same left-to-right column order they were
# authored in (this also controls the order new entries land in the shared
# string table: header then placeholder, column by column).
$fields = @(
    @{ Col = "Q";  Header = "Upcountry";                           Token = "{excel_data_line_item:upcountry}" },
    @{ Col = "R";  Header = "Upcountry Rate";                      Token = "{excel_data_line_item:upcountry_rate}" },
    @{ Col = "S";  Header = "Upcountry Max Distance Threshold";    Token = "{excel_data_line_item:upcountry_max_distance_threshold}" },
    @{ Col = "T";  Header = "Upcountry Approval";                  Token = "{excel_data_line_item:upcountry_approval}" },
    @{ Col = "U";  Header = "Upcountry Approval Email";            Token = "{excel_data_line_item:upcountry_approval_email}" },
    @{ Col = "V";  Header = "Invoice Email To";                    Token = "{excel_data_line_item:invoice_email_to}" },
    @{ Col = "W";  Header = "Invoice Email Cc";                    Token = "{excel_data_line_item:invoice_email_cc}" },
    @{ Col = "X";  Header = "Invoice Email Bcc";                   Token = "{excel_data_line_item:invoice_email_bcc}" },
    @{ Col = "Y";  Header = "PrePaid or Postpaid";                 Token = "{excel_data_line_item:is_prepaid}" },
    @{ Col = "Z";  Header = "PrePaid Amoun";                       Token = "{excel_data_line_item:prepaid_amount_limit}" },
    @{ Col = "AA"; Header = "PrePaid Notification Amount";         Token = "{excel_data_line_item:prepaid_notification_amount}" },
    @{ Col = "AB"; Header = "PostPaid Credit Period";              Token = "{excel_data_line_item:postpaid_credit_period}" },
    @{ Col = "AC"; Header = "PostPaid Notification Limit";         Token = "{excel_data_line_item:postpaid_notification_limit}" },
    @{ Col = "AD"; Header = "PostPaid Grace Period";               Token = "{excel_data_line_item:postpaid_grace_period}" }
)

foreach ($f in $fields) {
    $ws.Range($f.Col + "1").Value = $f.Header
    $ws.Range($f.Col + "2").Value = $f.Token
}

# Row 1 (header row) formatting: match the bold / left-aligned style already
# used across the existing header cells.
$ws.Range("N1").Copy() | Out-Null
$ws.Range("Q1:U1").PasteSpecial(-4122) | Out-Null
$ws.Range("W1:AB1").PasteSpecial(-4122) | Out-Null

$ws.Range("O1").Copy() | Out-Null
$ws.Range("V1").PasteSpecial(-4122) | Out-Null
$ws.Range("AC1:AD1").PasteSpecial(-4122) | Out-Null

# Row 2 (placeholder row) formatting: match the existing O2/P2 style.
$ws.Range("O2").Copy() | Out-Null
$ws.Range("Q2:AD2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Widen columns for the new, longer field names.
$ws.Columns.Item(16).ColumnWidth = 37.5   # P
$ws.Columns.Item(17).ColumnWidth = 28     # Q
$ws.Columns.Item(18).ColumnWidth = 30.5   # R
$ws.Columns.Item(19).ColumnWidth = 32     # S

# Print the sheet in portrait orientation.
$ws.PageSetup.Orientation = 1

# Restore the cursor/selection left behind by the last save.
$ws.Range("H22").Select() | Out-Null
